# The "Test Data" column of the two "testNonpositive" / negative-size rows in
# the importCandidatesHeader test-log table each contained a "header: ..."
# paragraph followed by a "line: 2" paragraph. Per the commit message ("...to
# either not include line numbers as test data..."), remove the paragraphs
# that merely report "line: <number>" as test data, leaving the "header: ..."
# paragraph intact.
#
# Walk the document's paragraphs back-to-front (so deleting one doesn't
# invalidate the indices of paragraphs we still need to visit) and delete any
# whose text is exactly "line: <digits>" (the paragraph mark is included in
# Range.Text, which is why we match with a wildcard rather than -eq).

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text
    if ($txt -like "line: *") {
        $p.Range.Delete()
    }
}
